$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.211.56"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.857.68"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'0.6991"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'237.16"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07667"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "'0.3044"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'0.08154"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.850.65"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'0.7152"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'5.146"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "'89.32"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "29.221.89"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "'5.746"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'13.23"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "'237.61"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'0.000007693"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "2.112.79"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'7.440"
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").Value = "'0.1471"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'162.19"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "'8.990"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "'2.001"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("D30").Value = "'1.419"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "'4.431"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "'1.482"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "'4.006"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").Value = "'0.05189"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "'1.162"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'0.9988"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "'2.659"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'0.01850"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "'2.720"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").Value = "'0.9320"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "1.146.63"
$ws.Range("E42").Value = "  +10.03%  "
$ws.Range("D43").Value = "'0.4280"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'70.76"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "'5.854"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "'1.793"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "2.009.89"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "'6.956"
$ws.Range("E51").Value = "  -3.75%  "
